$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A42 timestamp precision (same double value, re-saved by source tool)
$ws.Range("A42").Value = 44355.78620761227

# Add new row 43
$ws.Range("A43").Value = 44356.77938151157
$ws.Range("A43").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("B43").Value = 76251
$ws.Range("C43").Value = 64180
$ws.Range("D43").Value = 3455
$ws.Range("E43").Value = 2107
$ws.Range("F43").Value = 1482
$ws.Range("G43").Value = 20064
$ws.Range("H43").Value = 1460
$ws.Range("I43").Value = 884
$ws.Range("J43").Value = 201
